$d = $word.ActiveDocument

function Append-NewRun($searchText, $appendText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $end = $rng.End
    $newRng = $d.Range($end, $end)
    $newRng.InsertAfter($appendText)
    # Force a run split (new <w:r>) at the boundary without changing visible
    # formatting: toggle Bold on/off so the appended text becomes its own run
    # instead of being silently merged into the preceding run.
    $newRng.SetRange($end, $end + $appendText.Length)
    $newRng.Font.Bold = 1
    $newRng.Font.Bold = 0
}

# 1) rs6983267 location line: add the base-pair delta.
Append-NewRun "8:127401060 (GRCh38) --- 8:128413305 (GRCh37)" " = 1,012,245 "

# 2) rs4713266 location line: add the base-pair delta.
Append-NewRun "6:11218797 (GRCh38) --- 6:11219030 (GRCh37)" " = 233"

# 3) rs72699833 location line: add the base-pair delta.
Append-NewRun "1:121539689 (GRCh38) --- 1:121281487 (GRCh37)" " = -258,202"

# 4) Collapse the split "(1|1)" runs into a single contiguous run of text.
$startRng = $d.Content
$startRng.Find.Execute("homozygous dominate ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitStart = $startRng.End

$endRng = $d.Content
$endRng.Find.Execute("1), and heterozygous (0|1, 1|0) for the cancer-risk SNPs. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitEnd = $endRng.End

$mergeRng = $d.Range($splitStart, $splitEnd)
$mergeRng.Delete()
$insertRng = $d.Range($splitStart, $splitStart)
$insertRng.InsertAfter("(1|1), and heterozygous (0|1, 1|0) for the cancer-risk SNPs. ")

Write-Output "done"
